$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(119).Insert()

$ws.Cells.Item(119, 1).Value = 3
$ws.Cells.Item(119, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44518
$ws.Cells.Item(119, 5).Value = 5
$ws.Cells.Item(119, 6).Value = 100112001
$ws.Cells.Item(119, 7).Value = "Berenjena"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 115
$ws.Cells.Item(119, 11).Value = 7500
$ws.Cells.Item(119, 12).Value = 8000
$ws.Cells.Item(119, 13).Value = 7739
$ws.Cells.Item(119, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(119, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119, 16).Value = 129
$ws.Cells.Item(119, 17).Value = 60
$ws.Cells.Item(119, 18).Value = "Hortaliza"
